$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 14:44"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6291627
$ws.Range("C4").Value = 890
$ws.Range("E4").Value = 2553687

# Row 13 - Argentina
$ws.Range("D13").Value = 322461
$ws.Range("E13").Value = 107556
$ws.Range("G13").Value = 37
$ws.Range("H13").Value = 9155

# Row 18 - Arabia Saudita
$ws.Range("B18").Value = 318319
$ws.Range("C18").Value = 833
$ws.Range("D18").Value = 293964
$ws.Range("E18").Value = 20373
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = 3982

# Row 30 - Catar
$ws.Range("B30").Value = 119420
$ws.Range("C30").Value = 214
$ws.Range("D30").Value = 116313
$ws.Range("E30").Value = 2906
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 201

# Row 42 - Suecia
$ws.Range("B42").Value = 84729
$ws.Range("G42").Value = 5
$ws.Range("H42").Value = 5832

# Rows 44-46: "Paises Bajos" overtakes "Bielorrusia" and "Emiratos Arabes Unidos"
# in the ranking, so the three rows shift: row44 becomes Paises Bajos (updated
# figures), row45 becomes Bielorrusia (previous row44 figures), row46 becomes
# Emiratos Arabes Unidos (previous row45 figures).

# Row 46 used to be Paises Bajos -> becomes Emiratos Arabes Unidos (old row45 data)
$ws.Range("A46").Value = "Emiratos Arabes Unidos"
$ws.Range("B46").Value = 72154
$ws.Range("C46").Value = 614
$ws.Range("D46").Value = 62668
$ws.Range("E46").Value = 9099
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 387

# Row 45 used to be Emiratos Arabes Unidos -> becomes Bielorrusia (old row44 data)
$ws.Range("A45").Value = "Bielorrusia"
$ws.Range("B45").Value = 72302
$ws.Range("C45").Value = 161
$ws.Range("D45").Value = 71205
$ws.Range("E45").Value = 401
$ws.Range("G45").Value = 5
$ws.Range("H45").Value = 696

# Row 44 used to be Bielorrusia -> becomes Paises Bajos (new, updated figures)
$ws.Range("A44").Value = "Paises Bajos"
$ws.Range("B44").Value = 72464
$ws.Range("C44").Value = 601
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 6235

# Row 75 - Estado de Palestina
$ws.Range("B75").Value = 24471
$ws.Range("C75").Value = 596
$ws.Range("D75").Value = 16095
$ws.Range("E75").Value = 8209
$ws.Range("G75").Value = 5
$ws.Range("H75").Value = 167

# Row 82 - Dinamarca
$ws.Range("B82").Value = 17374
$ws.Range("C82").Value = 179
$ws.Range("D82").Value = 15499
$ws.Range("E82").Value = 1249

# Row 143 - Islandia
$ws.Range("B143").Value = 2128
$ws.Range("C143").Value = 7
$ws.Range("D143").Value = 2023

# Row 195 - Antigua y Barbuda
$ws.Range("B195").Value = 95
$ws.Range("C195").Value = 1
$ws.Range("E195").Value = 1
